$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.16"
$ws.Range("E2").Value = "'1.12%"
$ws.Range("D3").Value = "'36.12"
$ws.Range("E3").Value = "'-4.12%"
$ws.Range("D4").Value = "'5.077"
$ws.Range("E4").Value = "'1.37%"
$ws.Range("D5").Value = "'0.07853"
$ws.Range("E5").Value = "'-0.21%"
$ws.Range("D6").Value = "'2.200"
$ws.Range("E6").Value = "'-0.13%"
$ws.Range("D7").Value = "'7.925"
$ws.Range("E7").Value = "'-1.04%"
$ws.Range("D8").Value = "'0.9202"
$ws.Range("E8").Value = "'1.21%"
$ws.Range("D9").Value = "'0.09592"
$ws.Range("E9").Value = "'4.03%"
$ws.Range("D10").Value = "'0.1865"
$ws.Range("E10").Value = "'-0.56%"
$ws.Range("D11").Value = "'0.08655"
$ws.Range("E11").Value = "'2.31%"
$ws.Range("D12").Value = "'0.03478"
$ws.Range("E12").Value = "'-1.20%"
$ws.Range("D13").Value = "'0.09936"
$ws.Range("E13").Value = "'-0.01%"
$ws.Range("D14").Value = "'0.001430"
$ws.Range("E14").Value = "'-3.20%"
$ws.Range("D15").Value = "'0.005629"
$ws.Range("E15").Value = "'-0.14%"
$ws.Range("D16").Value = "'3.465"
$ws.Range("E16").Value = "'-0.35%"
$ws.Range("D17").Value = "'4.097"
$ws.Range("E17").Value = "'2.31%"
$ws.Range("D18").Value = "'2.458"
$ws.Range("E18").Value = "'19.33%"
$ws.Range("E19").Value = "'-1.08%"
$ws.Range("D20").Value = "'0.1301"
$ws.Range("E20").Value = "'0.12%"
$ws.Range("D21").Value = "'4.842"
$ws.Range("E21").Value = "'5.95%"
$ws.Range("E22").Value = "'-0.99%"
$ws.Range("D23").Value = "'0.04550"
$ws.Range("E23").Value = "'-1.99%"
$ws.Range("E24").Value = "'14.53%"
$ws.Range("D25").Value = "'0.001235"
$ws.Range("E25").Value = "'0.54%"
$ws.Range("E26").Value = "'8.06%"
$ws.Range("D27").Value = "'0.0004754"
$ws.Range("E27").Value = "'0.18%"
$ws.Range("E39").Value = "'4.71%"
$ws.Range("D40").Value = "'0.04778"
$ws.Range("E40").Value = "'1.52%"
$ws.Range("D41").Value = "'0.007723"
$ws.Range("E41").Value = "'-2.17%"
$ws.Range("E42").Value = "'0.61%"
$ws.Range("D43").Value = "'0.007729"
$ws.Range("E43").Value = "'0.90%"
$ws.Range("D44").Value = "'0.002233"
$ws.Range("E44").Value = "'-0.58%"
$ws.Range("D45").Value = "'0.01123"
$ws.Range("E45").Value = "'7.20%"
$ws.Range("D46").Value = "'0.00006260"
$ws.Range("E46").Value = "'0.16%"
$ws.Range("E47").Value = "'0.23%"
$ws.Range("D48").Value = "'0.0005805"
$ws.Range("E48").Value = "'0.08%"
$ws.Range("D49").Value = "'24.38"
$ws.Range("E49").Value = "'181.20%"
$ws.Range("D50").Value = "'0.002001"
$ws.Range("E50").Value = "'-25.83%"
$ws.Range("E51").Value = "'0.23%"
